# Fruta / hortaliza, semanal
# Rotates the weekly price-by-quality rows (2-23): each destination row
# takes on the Fecha/Calidad/Volumen/Precio* values that another row held,
# per the mapping below (rows 24-26 are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (values to copy from, read from the
# ORIGINAL/unmodified sheet state)
$map = @{
  2  = 5
  3  = 6
  4  = 7
  5  = 21
  6  = 22
  7  = 23
  8  = 12
  9  = 13
  10 = 14
  11 = 15
  12 = 16
  13 = 17
  14 = 18
  15 = 8
  16 = 9
  17 = 10
  18 = 11
  19 = 2
  20 = 3
  21 = 4
  22 = 19
  23 = 20
}

# Columns touched by the rotation (A1-style letters -> column index)
# D=4 (Fecha), L=12 (Calidad), M=13 (Volumen), N=14 (Precio minimo),
# O=15 (Precio maximo), P=16 (Precio promedio ponderado), S=19 (Precio $/Kg)
$cols = @(4, 12, 13, 14, 15, 16, 19)

# --- Phase 1: snapshot every original value we might need as a source ---
$snapshot = @{}
foreach ($r in $map.Keys) {
    foreach ($c in $cols) {
        $key = "$r,$c"
        $snapshot[$key] = $ws.Cells.Item($r, $c).Value2
    }
}

# --- Phase 2: write the rotated values from the snapshot ---
foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    foreach ($c in $cols) {
        $srcKey = "$srcRow,$c"
        $ws.Cells.Item($destRow, $c).Value2 = $snapshot[$srcKey]
    }
}
